$d = $word.ActiveDocument

function Split-RunAt($doc, $pos) {
    # Force a clean run boundary at an absolute character position by
    # inserting a paragraph break there and immediately deleting the
    # inserted paragraph mark again. Net effect on the document text is
    # nil, but Word is left with two separate runs instead of one.
    $r = $doc.Range($pos, $pos)
    $r.InsertParagraphAfter()
    $markEnd = $pos + 1
    $mark = $doc.Range($pos, $markEnd)
    $mark.Delete()
}

function Insert-TextAt($doc, $pos, $text) {
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
}

# ---------------------------------------------------------------------
# 1) cidades( codigo_ibge -> id ,nome,uf)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("codigo_ibge", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "id", 2)

# ---------------------------------------------------------------------
# 2) pessoas(...) : "logradouro" -> "endereco" (own run) and insert a
#    new "tipoUsuario," run before "cidade_id)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("logradouro", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "endereco", 2)

$rngStatus = $d.Content
$rngStatus.Find.Execute("obs,status,", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$posAfterStatus = $rngStatus.End
Insert-TextAt $d $posAfterStatus "tipoUsuario,"

# ---------------------------------------------------------------------
# 3) funcionários( -> funcionarios( split as "funciona" + "rios"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("funcionários", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "funcionarios", 2)

$rngFunc = $d.Content
$rngFunc.Find.Execute("funcionarios(", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$funcStart = $rngFunc.Start
$funcSplit = $funcStart + 8
Split-RunAt $d $funcSplit

# ---------------------------------------------------------------------
# 4) medicamentos(...) : insert ",unidade" run just before the closing ")"
# ---------------------------------------------------------------------
$rngMed = $d.Content
$rngMed.Find.Execute(",classe_terapeutica,tarja,posologia)", $true, $false, $false, `
                      $false, $false, $true, 1, $false, "", 0)
$medEnd = $rngMed.End
$medInsertPos = $medEnd - 1
Insert-TextAt $d $medInsertPos ",unidade"

# ---------------------------------------------------------------------
# 5) imagens(...) : insert ",descricao" run just before the closing ")"
# ---------------------------------------------------------------------
$rngImg = $d.Content
$rngImg.Find.Execute("paciente_id", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$rngImg2 = $d.Content
$rngImg2.Start = $rngImg.End
$rngImg2.End = $d.Content.End
$rngImg2.Find.Execute("link)", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$imgEnd = $rngImg2.End
$imgInsertPos = $imgEnd - 1
Insert-TextAt $d $imgInsertPos ",descricao"

# ---------------------------------------------------------------------
# 6) itensReceita(...) : ",dose,obs)" -> ",dose,obs" + ",periodo" +
#    bookmark "_GoBack" + ")"   (bookmark relocated here from procedimentos)
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$rngRec = $d.Content
$rngRec.Find.Execute(",dose,obs)", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$recEnd = $rngRec.End
$recInsertPos = $recEnd - 1
Insert-TextAt $d $recInsertPos ",periodo"

$rngRec2 = $d.Content
$rngRec2.Find.Execute(",periodo", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$bmPos = $rngRec2.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
